# Update Bhutan MSME summary figures:
#   - Enterprises density (per 1000 people): Micro 28.6 -> 28.59, SMEs 4.4 -> 4.39, MSMEs 33 -> 32.98
#   - Enterprises (% of total): SMEs 13.2 -> 13.22, MSMEs 99.4 -> 99.43
#
# The source values are stored as text (not numbers) in the workbook. Assigning a
# numeric-looking string straight to .Value would make Excel coerce the cell to a
# real number, so each value is entered with a leading apostrophe to force text,
# then the cell's Style is re-copied from an untouched neighbour in the same
# column (which still carries the original, non quote-prefixed style) so the
# cell formatting ends up identical to how it started.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11: "Enterprises density (per 1000 people)" -> Micro / SMEs / MSMEs
$ws.Range("B11").Value = "'28.59"
$ws.Range("B11").Style = $ws.Range("B10").Style

$ws.Range("C11").Value = "'4.39"
$ws.Range("C11").Style = $ws.Range("C10").Style

$ws.Range("D11").Value = "'32.98"
$ws.Range("D11").Style = $ws.Range("D10").Style

# Row 12: "Enterprises (% of total)" -> SMEs / MSMEs (Micro's 86.2 is unchanged)
$ws.Range("C12").Value = "'13.22"
$ws.Range("C12").Style = $ws.Range("C10").Style

$ws.Range("D12").Value = "'99.43"
$ws.Range("D12").Style = $ws.Range("D10").Style
